$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi following Dr Hou advice: recompute Lif-Lifr L-R pairs including the
# ECs sending cluster, adding 3 new rows (ECs -> ECs/FAPs/sCs) and refreshing
# all previously computed statistics for the existing FAPs rows.
$rows = @(
    @("ECs", "Lif", "Lifr", "ECs", 2, 0.6666666666666666, 0.904627, 2.713881, 0.2670571014571191, 0.2670571014571191, 3, 1, 14.96835266666667, 44.905058, 0.1240053612000741, 0.1240053612000741, 13.54077596778867, 121.866983710098, 0.03311651232723489, 0.03311651232723489),
    @("ECs", "Lif", "Lifr", "FAPs", 2, 0.6666666666666666, 0.904627, 2.713881, 0.2670571014571191, 0.2670571014571191, 3, 1, 80.77474466666666, 242.324234, 0.6691786071115035, 0.6691786071115035, 73.07101494357266, 657.639134492154, 0.1787088991723104, 0.1787088991723104),
    @("ECs", "Lif", "Lifr", "sCs", 2, 0.6666666666666666, 0.904627, 2.713881, 0.2670571014571191, 0.2670571014571191, 3, 1, 24.96420533333334, 74.892616, 0.2068160316884225, 0.2068160316884225, 22.58329417807733, 203.249647602696, 0.05523168995757381, 0.05523168995757381),
    @("FAPs", "Lif", "Lifr", "ECs", 3, 1, 2.482764666666667, 7.448294000000001, 0.7329428985428809, 0.7329428985428807, 3, 1, 14.96835266666667, 44.905058, 0.1240053612000741, 0.1240053612000741, 37.16289711900578, 334.466074071052, 0.0908888488728392, 0.09088884887283918),
    @("FAPs", "Lif", "Lifr", "FAPs", 3, 1, 2.482764666666667, 7.448294000000001, 0.7329428985428809, 0.7329428985428807, 3, 1, 80.77474466666666, 242.324234, 0.6691786071115035, 0.6691786071115035, 200.5446820174218, 1804.902138156796, 0.490469707939193, 0.4904697079391929),
    @("FAPs", "Lif", "Lifr", "sCs", 3, 1, 2.482764666666667, 7.448294000000001, 0.7329428985428809, 0.7329428985428807, 3, 1, 24.96420533333334, 74.892616, 0.2068160316884225, 0.2068160316884225, 61.98024693301156, 557.822222397104, 0.1515843417308487, 0.1515843417308487)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
